# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1885
    $ws.Range("F3").Value = 355
    $ws.Range("F4").Value = 1150
    $ws.Range("F5").Value = 1174
    $ws.Range("F7").Value = 5979
    $ws.Range("F8").Value = 99
}
